$wb = $excel.ActiveWorkbook

# Add a new sheet (becomes the new active/last sheet) and rename it
$ws = $wb.Worksheets.Add()
$ws.Name = "ValidLogin"

# Set new cell values: un/pw headers, admin/manager values
$ws.Range("A1").Value = "un"
$ws.Range("B1").Value = "pw"
$ws.Range("A2").Value = "admin"
$ws.Range("B2").Value = "manager"

[void]$ws.Range("B2").Select()

# Remove the original default sheet
[void]$wb.Worksheets.Item("Sheet1").Delete()
